# #5: cash & deposit done
# Adds the standard metadata columns (property_category/category/date/
# legislator_name/legislator_id/source_file/index) to the 現金 (cash) and
# 存款 (deposit) sheets, matching the layout already used by 土地/建物/汽車,
# and turns row 1 of each sheet into a proper header row.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-TextValue($range, $text) {
    # Prefix with an apostrophe so date-shaped strings like "2012-03-01"
    # are stored as literal text instead of being auto-parsed into a date
    # serial number.
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------------
# Sheet 4 - 現金 (cash)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Row 1 becomes a real header row: B1 keeps its "currency" meaning while
# C1/D1 switch from data values to header labels.
$ws4.Range("B1").Value = "currency"
$ws4.Range("C1").Value = "owner"
$ws4.Range("D1").Value = "total"

$ws4Headers = @("property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
$ws4HeaderCols = @("E", "F", "G", "H", "I", "J", "K")
for ($i = 0; $i -lt $ws4Headers.Length; $i++) {
    $cell = $ws4.Range("$($ws4HeaderCols[$i])1")
    Set-TextValue $cell $ws4Headers[$i]
    $ws4.Range("B1").Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
}

# Data rows 2-5: append the same seven metadata columns. "index" (column K)
# repeats the row's own index value (column A).
for ($r = 2; $r -le 5; $r++) {
    $idxValue = $ws4.Range("A$r").Value2

    Set-TextValue $ws4.Range("E$r") "cash"
    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("E$r").PasteSpecial($xlPasteFormats) | Out-Null

    Set-TextValue $ws4.Range("F$r") "normal"
    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("F$r").PasteSpecial($xlPasteFormats) | Out-Null

    Set-TextValue $ws4.Range("G$r") "2012-03-01"
    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("G$r").PasteSpecial($xlPasteFormats) | Out-Null

    Set-TextValue $ws4.Range("H$r") "林淑芬"
    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("H$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws4.Range("I$r").Value = 1337
    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("I$r").PasteSpecial($xlPasteFormats) | Out-Null
    $ws4.Range("I$r").Value = 1337

    Set-TextValue $ws4.Range("J$r") "tmp3f851"
    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("J$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws4.Range("B$r").Copy() | Out-Null
    $ws4.Range("K$r").PasteSpecial($xlPasteFormats) | Out-Null
    $ws4.Range("K$r").Value = $idxValue
}

# ---------------------------------------------------------------------
# Sheet 5 - 存款 (deposit)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Row 1 becomes a real header row too.
$ws5.Range("B1").Value = "bank"
$ws5.Range("C1").Value = "deposit_type"
$ws5.Range("D1").Value = "currency"
$ws5.Range("E1").Value = "owner"
$ws5.Range("F1").Value = "total"

$ws5Headers = @("property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
$ws5HeaderCols = @("G", "H", "I", "J", "K", "L", "M")
for ($i = 0; $i -lt $ws5Headers.Length; $i++) {
    $cell = $ws5.Range("$($ws5HeaderCols[$i])1")
    Set-TextValue $cell $ws5Headers[$i]
    $ws5.Range("B1").Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
}

# Data rows 2-5: append the same seven metadata columns.
for ($r = 2; $r -le 5; $r++) {
    $idxValue = $ws5.Range("A$r").Value2

    Set-TextValue $ws5.Range("G$r") "deposit"
    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("G$r").PasteSpecial($xlPasteFormats) | Out-Null

    Set-TextValue $ws5.Range("H$r") "normal"
    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("H$r").PasteSpecial($xlPasteFormats) | Out-Null

    Set-TextValue $ws5.Range("I$r") "2012-03-01"
    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("I$r").PasteSpecial($xlPasteFormats) | Out-Null

    Set-TextValue $ws5.Range("J$r") "林淑芬"
    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("J$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws5.Range("K$r").Value = 1337
    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("K$r").PasteSpecial($xlPasteFormats) | Out-Null
    $ws5.Range("K$r").Value = 1337

    Set-TextValue $ws5.Range("L$r") "tmp3f851"
    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("L$r").PasteSpecial($xlPasteFormats) | Out-Null

    $ws5.Range("B$r").Copy() | Out-Null
    $ws5.Range("M$r").PasteSpecial($xlPasteFormats) | Out-Null
    $ws5.Range("M$r").Value = $idxValue
}

# F4 was previously stored as the text "1527000"; the corrected sheet
# stores it as a real number, same as the other amount cells.
$ws5.Range("F4").Value = 1527000
